$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: Label, URL, Picture
# (shared strings are added in the order values are assigned, so write the
# URL first, then the label, then the picture filename, to match the
# target shared-string ordering)
$ws.Range("B3").Value = "https://www.crunchyroll.com/series/GDKHZEJ0K/solo-leveling"
$ws.Range("A3").Value = "Solo Levelling"
$ws.Range("C3").Value = "SoloLeveling.jpg"

# Add hyperlink on the URL cell (B3), matching the style of existing B2 hyperlink
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.crunchyroll.com/series/GDKHZEJ0K/solo-leveling")
$ws.Range("B3").Style = "Hyperlink"

# Update the selection to match the target state
$ws.Range("E5").Select()
